$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) of the last existing row (row 20) into the new row 21
$ws.Range("A20:I20").Copy($ws.Range("A21:I21"))

# Fill in the new workout record for 2020-06-24 (serial date 44006)
$ws.Cells.Item(21, 1).Value = 44006       # DATETIME
$ws.Cells.Item(21, 2).Value = 100         # WAIST
$ws.Cells.Item(21, 3).Value = 106         # HIP
$ws.Cells.Item(21, 4).Value = 0.9         # WHR
$ws.Cells.Item(21, 5).Value = "SAME"      # WHR_IMPROVEMENT
$ws.Cells.Item(21, 6).Value = 0           # WHR_CHANGE
$ws.Cells.Item(21, 7).Value = 82.1        # WEIGHT
$ws.Cells.Item(21, 8).Value = 28.7        # BMI
$ws.Cells.Item(21, 9).Value = "OVERWEIGHT" # OBESITY

# Select the whole sheet (mirrors the selection state captured in the saved file)
$ws.Cells.Select()
